# Weekly fruit/vegetable price update: a new price record (week) is
# inserted into the "Haba" (broad bean) sheet at row 60, pushing the
# existing records (previously rows 60-99) down by one row to rows 61-100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60, shifting rows 60:99 down to 61:100.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new weekly record.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 45236
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100112026
$ws.Range("G60").Value = "Haba"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 220
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = 10909
$ws.Range("N60").Value = "$/saco 25 kilos"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 436
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"
